# Edit script: update column G ("K") values for rows 2-37 in Sheet1
# per commit: "regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$kValues = @{
    2 = 5
    3 = 7
    4 = 8
    5 = 4
    6 = 5
    7 = 14
    8 = 5
    9 = 5
    10 = 3
    11 = 6
    12 = 6
    13 = 9
    14 = 6
    15 = 11
    16 = 5
    17 = 10
    18 = 10
    19 = 5
    20 = 9
    21 = 10
    22 = 8
    23 = 8
    24 = 6
    25 = 11
    26 = 5
    27 = 5
    28 = 7
    29 = 7
    30 = 10
    31 = 11
    32 = 10
    33 = 8
    34 = 4
    35 = 6
    36 = 8
    37 = 5
}

foreach ($row in $kValues.Keys) {
    $ws.Cells.Item([int]$row, 7).Value = $kValues[$row]
}
